# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 unicode char used in one ShibaInu price cell (e.g. 0.0₃0703)
$sub3 = [char]0x2083
$d17val = [string]::Concat("0.0", $sub3, "0703")

# Cells in column D hold price text that can look like plain numbers
# (e.g. "0.474", "206.12"); force them to remain Text so Excel does not
# auto-convert them to numeric values.
$priceUpdates = @(
    @{Cell="D2"; Value='25.170.75'},
    @{Cell="D3"; Value='1.550.40'},
    @{Cell="D5"; Value='206.12'},
    @{Cell="D7"; Value='0.474'},
    @{Cell="D9"; Value='0.239'},
    @{Cell="D10"; Value='17.62'},
    @{Cell="D11"; Value='0.0779'},
    @{Cell="D12"; Value='1.768.33'},
    @{Cell="D13"; Value='1.545.27'},
    @{Cell="D14"; Value='3.95'},
    @{Cell="D16"; Value='25.168.83'},
    @{Cell="D17"; Value=$d17val},
    @{Cell="D18"; Value='58.32'},
    @{Cell="D20"; Value='183.83'},
    @{Cell="D21"; Value='4.08'},
    @{Cell="D22"; Value='9.17'},
    @{Cell="D23"; Value='5.80'},
    @{Cell="D25"; Value='139.54'},
    @{Cell="D26"; Value='0.127'},
    @{Cell="D30"; Value='1.15'},
    @{Cell="D31"; Value='0.0459'},
    @{Cell="D33"; Value='2.95'},
    @{Cell="D36"; Value='1.083.64'},
    @{Cell="D39"; Value='2.24'},
    @{Cell="D40"; Value='0.489'},
    @{Cell="D41"; Value='0.800'},
    @{Cell="D42"; Value='0.751'},
    @{Cell="D43"; Value='92.14'},
    @{Cell="D44"; Value='5.00'},
    @{Cell="D45"; Value='1.683.60'},
    @{Cell="D47"; Value='52.00'},
    @{Cell="D48"; Value='0.0502'}
)

# Cells in columns B, C, E are never ambiguous with numbers (names, URLs,
# or percentage strings padded with spaces), so they can be set directly.
$otherUpdates = @(
    @{Cell="E2"; Value='  -2.83%  '},
    @{Cell="E3"; Value='  -4.37%  '},
    @{Cell="E5"; Value='  -3.48%  '},
    @{Cell="E6"; Value='  -0.03%  '},
    @{Cell="E7"; Value='  -5.55%  '},
    @{Cell="E8"; Value='  -2.17%  '},
    @{Cell="E9"; Value='  -3.98%  '},
    @{Cell="E10"; Value='  -3.67%  '},
    @{Cell="E11"; Value='  -1.05%  '},
    @{Cell="E12"; Value='  -4.32%  '},
    @{Cell="E13"; Value='  -4.44%  '},
    @{Cell="E14"; Value='  -5.27%  '},
    @{Cell="E15"; Value='  -4.91%  '},
    @{Cell="E16"; Value='  -2.84%  '},
    @{Cell="E17"; Value='  -4.00%  '},
    @{Cell="E18"; Value='  -4.48%  '},
    @{Cell="E20"; Value='  -3.91%  '},
    @{Cell="E21"; Value='  -3.48%  '},
    @{Cell="E22"; Value='  -4.21%  '},
    @{Cell="E23"; Value='  -4.35%  '},
    @{Cell="E24"; Value='  -0.04%  '},
    @{Cell="B25"; Value='Monero'},
    @{Cell="C25"; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Cell="E25"; Value='  -2.80%  '},
    @{Cell="B26"; Value='Stellar'},
    @{Cell="C26"; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'},
    @{Cell="E26"; Value='  -4.23%  '},
    @{Cell="E27"; Value='  -3.89%  '},
    @{Cell="E28"; Value='  -2.68%  '},
    @{Cell="E29"; Value='  -5.20%  '},
    @{Cell="E30"; Value='  -6.73%  '},
    @{Cell="E31"; Value='  -4.64%  '},
    @{Cell="E32"; Value='  -3.21%  '},
    @{Cell="E33"; Value='  -4.68%  '},
    @{Cell="E34"; Value='  -3.72%  '},
    @{Cell="E35"; Value='  -4.02%  '},
    @{Cell="E36"; Value='  -3.12%  '},
    @{Cell="E37"; Value='  -0.15%  '},
    @{Cell="E38"; Value='  -2.99%  '},
    @{Cell="B39"; Value='MXToken'},
    @{Cell="C39"; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Cell="E39"; Value='  -7.66%  '},
    @{Cell="B40"; Value='ImmutableX'},
    @{Cell="C40"; Value='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'},
    @{Cell="E40"; Value='  -5.47%  '},
    @{Cell="E41"; Value='  +4.67%  '},
    @{Cell="E42"; Value='  -11.16%  '},
    @{Cell="E43"; Value='  -5.93%  '},
    @{Cell="E44"; Value='  -2.74%  '},
    @{Cell="E45"; Value='  -4.26%  '},
    @{Cell="E46"; Value='  -7.09%  '},
    @{Cell="E47"; Value='  -4.13%  '},
    @{Cell="E48"; Value='  -5.26%  '},
    @{Cell="E49"; Value='  -2.71%  '},
    @{Cell="E50"; Value='  -1.64%  '},
    @{Cell="E51"; Value='  -0.28%  '}
)

foreach ($u in $priceUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

foreach ($u in $otherUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
